$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.768.38'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").Value = '2.493.67'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '535.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.566'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.64%  '
$ws.Range("D9").Value = '2.518.06'
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("E10").Value = '  +0.03%  '
$ws.Range("E11").Value = '  -2.77%  '
$ws.Range("E12").Value = '  -2.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.347'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.89%  '
$ws.Range("D14").Value = '2.942.47'
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.90'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.38%  '
$ws.Range("D16").Value = '58.701.22'
$ws.Range("E16").Value = '  -0.94%  '
$ws.Range("E17").Value = '  -1.39%  '
$ws.Range("D18").Value = '2.520.06'
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.90'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.92'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.419'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.164'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.51%  '
$ws.Range("D29").Value = '0.0₃0767'
$ws.Range("E29").Value = '  -1.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.62'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '171.04'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.50%  '
$ws.Range("E32").Value = '  -2.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.16'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.38%  '
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.40'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.34'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.07'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.53'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.66'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.806'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '283.80'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.16'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.995'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.607'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '130.61'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.93%  '
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0921'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0505'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.09%  '
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.30'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.58%  '
